$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-79
$iVals = @(
    8, 8, 9, 8, 8, 8, 8, 8, 8, 8, 8, 9, 7, 8, 8, 8, 8, 8, 8, 7, 8, 8, 8, 7, 8, 7, 10, 7, 7, 8, 8, 8, 8, 8, 8, 8, 10, 8, 7, 8, 7, 8, 11, 8, 7, 8, 8, 9, 7, 8, 8, 9, 8, 9, 8, 8, 7, 10, 8, 8, 7, 7, 8, 7, 7, 10, 9, 7, 10, 8, 7, 5, 4, 6, 8, 6, 4, 5
)
$jVals = @(
    8, 8, 9, 8, 8, 8, 8, 8, 8, 8, 8, 9, 7, 8, 8, 8, 8, 8, 8, 7, 8, 8, 8, 7, 8, 8, 10, 8, 7, 8, 8, 8, 8, 8, 8, 8, 10, 8, 7, 8, 8, 8, 11, 8, 7, 8, 8, 9, 8, 8, 8, 9, 8, 9, 8, 8, 8, 10, 8, 8, 7, 7, 8, 7, 7, 10, 9, 7, 10, 8, 7, 5, 4, 7, 8, 6, 4, 5
)

for ($r = 2; $r -le 79; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
